# Add a new applicant row (row 24) to the bottom of the Qabul list.
# Mirrors the formatting of the existing data rows (3-23): plain text
# cells with no explicit style, even though some values look numeric
# or date-like (JSHIR id, phone numbers, the "Sana" date column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A24:K24")

# Force text interpretation first so values such as the 14-digit JSHIR
# code, the "+998..." phone numbers and the "2025-05-07" date string are
# not auto-converted into numbers / a date serial by Excel's input
# parsing.
$newRow.NumberFormat = "@"

$ws.Range("A24").Value = "Amirov Akrom"
$ws.Range("B24").Value = "Yurisprudensiya"
$ws.Range("C24").Value = "O'zbek tili"
$ws.Range("D24").Value = "Kunduzgi"
$ws.Range("E24").Value = "AA1234567"
$ws.Range("F24").Value = "12345678901234"
$ws.Range("G24").Value = "Toshkent viloyati"
$ws.Range("H24").Value = "Olmaliq shahri"
$ws.Range("I24").Value = "+998939849910"
$ws.Range("J24").Value = "+998945289910"
$ws.Range("K24").Value = "2025-05-07"

# Drop the temporary "@" text format again so the new row ends up with
# the same unstyled look as the other data rows (no explicit cell
# style), now that the values are safely stored as text.
$newRow.ClearFormats()
